# Automatische test-sync: 2025-08-19 20:54:50
# Append the new mail-log entry (row 19) to the "Logs" sheet, extend the
# conditional-formatting ranges that tracked the previous last row (18),
# and bump the "Dashboard" summary count for the affected category.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the Logs sheet ---------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(19, 1).Value = "Nieuwe bestelling"
$logs.Cells.Item(19, 2).Value = "planning@testbedrijf123.nl"
$logs.Cells.Item(19, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(19, 6).Value = "2025-08-19 20:54:06"
$logs.Cells.Item(19, 7).Value = "Nee"
$logs.Cells.Item(19, 8).Value = "Ja"
$logs.Cells.Item(19, 9).Value = "Nee"
$logs.Cells.Item(19, 10).Value = "Nee"

# --- 2. Extend the conditional formatting ranges from row 18 to row 19 --
$colsToExtend = "D", "G", "H", "I", "J"
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "18")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "19")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the Dashboard rollup count for the category -------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 18
